$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 89
$ws.Range("A89").Value = "'-673"
$ws.Range("B89").Value = "'11/17/2025"
$ws.Range("C89").Value = "Baunes 2195"
$ws.Range("D89").Value = "'12"
$ws.Range("E89").Value = "Pendiente ADM"
$ws.Range("F89").Value = "NEW"
$ws.Range("G89").Value = "Pendiente"
$ws.Range("H89").Value = "Sacar columna con prioridad  Tensar linga en columna continua"
$ws.Range("I89").Value = 1
$ws.Range("J89").Value = "Desmonte"
$ws.Range("K89").Value = "Sin equipos"
$ws.Range("L89").Value = "Pasante"
$ws.Range("M89").Value = -58.485702
$ws.Range("N89").Value = -34.576702
$ws.Range("O89").Value = "Paternal"
$ws.Range("P89").Value = "Capital Norte"
$ws.Range("Q89").Value = "ATH-C"
$ws.Range("R89").Value = "Fuera de Poligono OVL"

# Row 90
$ws.Range("A90").Value = "'-674"
$ws.Range("B90").Value = "'11/17/2025"
$ws.Range("C90").Value = "Miller 3597"
$ws.Range("D90").Value = "'12"
$ws.Range("E90").Value = "Pendiente ADM"
$ws.Range("F90").Value = "NEW"
$ws.Range("G90").Value = "Pendiente"
$ws.Range("H90").Value = "Aplomar/ enderezar columna y colocar rienda a pique"
$ws.Range("I90").Value = 1
$ws.Range("J90").Value = "Aplomo"
$ws.Range("K90").Value = "Sin equipos"
$ws.Range("L90").Value = "Pasante"
$ws.Range("M90").Value = -58.489898
$ws.Range("N90").Value = -34.560497
$ws.Range("O90").Value = "Saavedra"
$ws.Range("P90").Value = "Capital Norte"
$ws.Range("Q90").Value = "PUE-F"
$ws.Range("R90").Value = "Fuera de Poligono OVL"
